
$d = $word.ActiveDocument

# --- 1) Heading paragraph: new title + line break + new arxiv URL ---
$d.Paragraphs(1).Range.Text = 'Review 160: [Short] Vision Transformers Need Registers' + [char]11 + 'https://arxiv.org/abs/2309.16588'

# --- 2) Bold "Paper:" paragraph: new arxiv URL (v2) ---
$d.Paragraphs(2).Range.Text = 'Paper: https://arxiv.org/abs/2309.16588v2'

# --- 3) Remove the now-redundant standalone URL paragraph (old paragraph 4) ---
$urlPara = $d.Paragraphs(4)
$urlRange = $d.Range($urlPara.Range.Start, $urlPara.Range.End)
$urlRange.Delete()

# --- 4) Rewrite the four Hebrew body paragraphs (now at indices 6,8,10,12 after the delete) ---
$d.Paragraphs(6).Range.Text = 'ממש אהבתי את המאמר זה אבל לא בגלל שהוא הציע איזה רעיון מהפכני (הרעיון די נחמד אבל לא איזו פצצה). הסיבה לכך היא שהמאמר הזה הדגיש שוב את העובדה כמה מעט אנחנו מבינים מודלי ענק בעל ביליוני פרמטרים'
$d.Paragraphs(8).Range.Text = 'היום ב- #shorthebrewpapereviews אנו סוקרים מאמר שבמילים ממש פשוטות מצא שייצוג הדאטה שהטרנספורמרים הוויזואליים(כמו DINOv2) מפיקים מכילים דאטה מיותר שלא תורם לביצועי המודל יותר מדי. אזכיר שהמטרה העיקרית של מודלי הענק האלו היא לבנות ייצוג של דאטה המכיל את הפיצ''רים המהותיים ביותר שלו. '
$d.Paragraphs(10).Range.Text = 'כלומר הטרנספורמרים הויזאליים לא מצליחים לקודד את המידע בצורה המיטבית יש חלקים מיותרים בייצוג הזה. איך המחברים בכלל הגיעו לזה? הם שמו לב שיש פאצ''ים בתמונה שנורמה של ייצוגם (מהשכבה האחרונה) היא גדולה באופן אנומלי יחסית לייצוגי הפאצ''ים האחרים.'
$d.Paragraphs(12).Range.Text = 'המחברים גם שמו לב שייצוגים של פאצ''ים חריגים אלו מאוד דומים לייצוגי הפאצ''ים הסמוכים (מבחינת מרחק קוסיין).  בנוסף יכולת של ייצוגי פאצ''ים אנומליים אלו להציג את מיקום הפאץ'' בתמונה היא משמעותית יותר נמוכה מהפאצ''ם הרגילים (אימנו מודל לזיהוי המיקום).  הם עשו עוד בדיקות נוספות ששיכנעו אותם שייצוג הפאצ''ים האלו לא משפר את איכות המודל. '

# --- 5) Insert a brand-new paragraph after the empty paragraph following body4 (now index 13) ---
$anchorPara = $d.Paragraphs(13)
$newPara = $anchorPara.Range.InsertParagraphAfter()
$d.Paragraphs(14).Range.Text = 'אז מה הם עשו? משהו די אלגנטי (זה לא רעיון חדש כי כבר עשו זאת לפני כמה שנים במאמר על מודלי שפה). אז הם הוסיפו כמה טוקנים (אחרי טוקן cls) שמטרתם היא להכיל מידע לא רלוונטי. ייצוגי טוקנים אלה פשוט נזרקים ולא משמשים לא לאימון ולא לאינפרנס. וזה אכן משפר את ביצועי המודל בכמה משימות.'

Write-Output "edit complete"
